# Scheduled runner update: refresh cached market-price / profit figures
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2515.7036
$ws.Range("I40").Value = 2675.2942
$ws.Range("J40").Value = 2244.4
$ws.Range("K40").Value = 2675.2942
$ws.Range("L40").Value = 2244.4
$ws.Range("M40").Value = -2500.2942
$ws.Range("N40").Value = -2594.4

$ws.Range("H54").Value = 2076
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H64").Value = 7701.6
$ws.Range("I64").Value = 7701.6
$ws.Range("K64").Value = 7701.6
$ws.Range("M64").Value = -7453.6

$ws.Range("H67").Value = 7701.6
$ws.Range("I67").Value = 7701.6
$ws.Range("K67").Value = 7701.6
$ws.Range("M67").Value = -6843.6

$ws.Range("H74").Value = 9857.429
$ws.Range("I74").Value = 10003
$ws.Range("K74").Value = 10003
$ws.Range("M74").Value = -9067

$ws.Range("H77").Value = 9857.429
$ws.Range("I77").Value = 10003
$ws.Range("K77").Value = 50015
$ws.Range("M77").Value = -45335

$ws.Range("H103").Value = 1689.6875
$ws.Range("I103").Value = 1595.3334
$ws.Range("J103").Value = 1811
$ws.Range("K103").Value = 4786.0002
$ws.Range("L103").Value = 5433
$ws.Range("M103").Value = -4200.0002
$ws.Range("N103").Value = -6605

$ws.Range("H132").Value = 1686.8572
$ws.Range("I132").Value = 1162
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 3486
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -956
$ws.Range("N132").Value = -14057

$ws.Range("H138").Value = 2662.2974
$ws.Range("J138").Value = 3450.9048
$ws.Range("L138").Value = 10352.7144
$ws.Range("N138").Value = -20632.7144

$ws.Range("H141").Value = 1073.3334
$ws.Range("I141").Value = 901.25
$ws.Range("J141").Value = 1417.5
$ws.Range("K141").Value = 2703.75
$ws.Range("L141").Value = 4252.5
$ws.Range("M141").Value = 2476.25
$ws.Range("N141").Value = -14612.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 368.7143
$ws.Range("I25").Value = 188.66667
$ws.Range("K25").Value = 188.66667
$ws.Range("M25").Value = 213.33333

$ws.Range("H35").Value = 5025.6665
$ws.Range("I35").Value = 1274
$ws.Range("K35").Value = 1274
$ws.Range("M35").Value = -868

$ws.Range("H55").Value = 31621.143
$ws.Range("J55").Value = 33960
$ws.Range("L55").Value = 33960
$ws.Range("N55").Value = -34590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 1602.4286
$ws.Range("I24").Value = 636.4
$ws.Range("J24").Value = 4017.5
$ws.Range("K24").Value = 636.4
$ws.Range("L24").Value = 4017.5
$ws.Range("M24").Value = -401.4
$ws.Range("N24").Value = -4487.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1873.6842
$ws.Range("I22").Value = 763.5
$ws.Range("K22").Value = 763.5
$ws.Range("M22").Value = -413.5

$ws.Range("H37").Value = 23500
$ws.Range("I37").Value = 20500
$ws.Range("K37").Value = 20500
$ws.Range("M37").Value = -20393

$ws.Range("H51").Value = 36823.75
$ws.Range("J51").Value = 34431.668
$ws.Range("L51").Value = 34431.668
$ws.Range("N51").Value = -35903.668

$ws.Range("H61").Value = 36823.75
$ws.Range("J61").Value = 34431.668
$ws.Range("L61").Value = 34431.668
$ws.Range("N61").Value = -35127.668

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 4061.459
$ws.Range("I134").Value = 3684.0566
$ws.Range("K134").Value = 11052.1698
$ws.Range("M134").Value = -8517.1698

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3982
$ws.Range("I3").Value = 3979
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 11937
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = -11825
$ws.Range("N3").Value = -12224

$ws.Range("H26").Value = 18091.486
$ws.Range("J26").Value = 68788.89
$ws.Range("L26").Value = 206366.67
$ws.Range("N26").Value = -206942.67

$ws.Range("H70").Value = 500
$ws.Range("I70").Value = 500
$ws.Range("K70").Value = 1500
$ws.Range("M70").Value = -1185

$ws.Range("H73").Value = 500
$ws.Range("I73").Value = 500
$ws.Range("K73").Value = 1500
$ws.Range("M73").Value = -408

$ws.Range("H75").Value = 4909
$ws.Range("I75").Value = 500
$ws.Range("K75").Value = 1500
$ws.Range("M75").Value = -502

$ws.Range("H78").Value = 4909
$ws.Range("I78").Value = 500
$ws.Range("K78").Value = 4500
$ws.Range("M78").Value = 492

$ws.Range("H114").Value = 4176.6
$ws.Range("I114").Value = 3970.75
$ws.Range("K114").Value = 11912.25
$ws.Range("M114").Value = -8658.25

$ws.Range("H131").Value = 17551848
$ws.Range("I131").Value = 111112010
$ws.Range("J131").Value = 9316.6875
$ws.Range("K131").Value = 333336030
$ws.Range("L131").Value = 27950.0625
$ws.Range("M131").Value = -333330990
$ws.Range("N131").Value = -38030.0625

$ws.Range("H133").Value = 2152.375
$ws.Range("I133").Value = 2152.375
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 6457.125
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6647.033
$ws.Range("I46").Value = 1991.2
$ws.Range("J46").Value = 8974.950000000001
$ws.Range("K46").Value = 1991.2
$ws.Range("L46").Value = 8974.950000000001
$ws.Range("M46").Value = -1803.2
$ws.Range("N46").Value = -9350.950000000001

$ws.Range("H132").Value = 9444.223
$ws.Range("I132").Value = 4266.6665
$ws.Range("K132").Value = 12799.9995
$ws.Range("M132").Value = -10269.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 56735.4
$ws.Range("J54").Value = 20919.25
$ws.Range("L54").Value = 20919.25
$ws.Range("N54").Value = -21959.25

$ws.Range("H113").Value = 725.3
$ws.Range("I113").Value = 821.7143
$ws.Range("J113").Value = 500.33334
$ws.Range("K113").Value = 2465.1429
$ws.Range("L113").Value = 1501.00002
$ws.Range("M113").Value = -295.1428999999998
$ws.Range("N113").Value = -5841.000019999999

$ws.Range("H132").Value = 4546.8623
$ws.Range("I132").Value = 2233.182
$ws.Range("J132").Value = 5960.778
$ws.Range("K132").Value = 6699.545999999999
$ws.Range("L132").Value = 17882.334
$ws.Range("M132").Value = -4169.545999999999
$ws.Range("N132").Value = -22942.334

